$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 295; this shifts the existing rows 295:406 down to
# 296:407 (carrying their values/styles with them), matching the target
# diff where each row n (296..406) ends up holding what used to be in
# row n-1, and a new row 407 is appended holding what used to be the last
# row (406).
$ws.Rows.Item(295).Insert()

# Populate the newly inserted row 295 with the new record's data.
$ws.Range("A295").Value = 5
$ws.Range("B295").Value = "Macroferia Regional de Talca"
$ws.Range("C295").Value = "Maule"
$ws.Range("D295").Value = 44795
$ws.Range("E295").Value = 7
$ws.Range("F295").Value = 100112032
$ws.Range("G295").Value = "Zapallo italiano"
$ws.Range("H295").Value = "Sin especificar"
$ws.Range("I295").Value = "Primera"
$ws.Range("J295").Value = 200
$ws.Range("K295").Value = 22000
$ws.Range("L295").Value = 22000
$ws.Range("M295").Value = 22000
$ws.Range("N295").Value = "$/caja 50 unidades"
$ws.Range("O295").Value = "Región de Arica y Parinacota"
$ws.Range("P295").Value = 440
$ws.Range("Q295").Value = 50
$ws.Range("R295").Value = "Hortaliza"
